$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Brand"
$ws.Range("B1").Value = "Model"
$ws.Range("C1").Value = "Images"

# Row 2: Sikander DI 35
$ws.Range("A2").Value = "Sonalika Tractors"
$ws.Range("B2").Value = "Sikander DI 35"
$ws.Range("C2").Value = "['SikanderDI35img0-35-di-sikander-1631525637.png', 'SikanderDI35img1-upload-1631525637-0.png', 'SikanderDI35img2-upload-1631525637-0.png']"

# Row 3: DI 734 (S1)
$ws.Range("A3").Value = "Sonalika Tractors"
$ws.Range("B3").Value = "DI 734 (S1)"
$ws.Range("C3").Value = "['DI734(S1)img0-di-734-s1-1631528609.png', 'DI734(S1)img1-upload-1631528609-0.png', 'DI734(S1)img2-di-734-s1-1631528609.png']"

# Row 4: DI 47 RX (new row)
$ws.Range("A4").Value = "Sonalika Tractors"
$ws.Range("B4").Value = "DI 47 RX"
$ws.Range("C4").Value = "['DI47RXimg0-di-47-rx-1631526733.png', 'DI47RXimg1-upload-1631526733-0.png', 'DI47RXimg2-upload-1631526733-0.png']"

# Row 5: Tiger DI 50 4WD (new row)
$ws.Range("A5").Value = "Sonalika Tractors"
$ws.Range("B5").Value = "Tiger DI 50 4WD"
$ws.Range("C5").Value = "['TigerDI504WDimg0-sonalika-tiger-di-50-4wd-1696503646.png', 'TigerDI504WDimg1-sonalika-tiger-di-50-4wd-16965036460.png', 'TigerDI504WDimg2-sonalika-tiger-di-50-4wd-1696503646.png']"
